# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" sheet right after "总计" (i.e. before "2022-Q2"),
#    and populate it with the fund-holdings table.
# 2) Prepend the matching summary row to the "总计" sheet, shifting the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q3" worksheet
# ---------------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row (same header text/style as the other quarterly sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q3.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Fund holdings data: B=code, C=name, D=size, E=stock position, F=weight,
# G=market value (亿元), H=position rank.
# G values are text like the rest of the numeric-looking columns, except
# where the source workbook stores an exact-zero market value as a real
# number (marked here with the "NUM:" prefix).
$fundRows = @(
    @('161725', '招商中证白酒指数A', '609.75', '94.82', '1.96', '11.9511', 9),
    @('012414', '招商中证白酒指数C', '92.42', '94.82', '1.96', '1.8114', 9),
    @('003378', '泰康策略优选灵活配置混合', '13.78', '83.10', '4.77', '0.6573', 4),
    @('010874', '泰康品质生活混合A', '6.86', '84.42', '5.02', '0.3444', 4),
    @('010875', '泰康品质生活混合C', '3.44', '84.42', '5.02', '0.1727', 4),
    @('005014', '泰康景泰回报混合A', '8.99', '34.29', '1.50', '0.1348', 5),
    @('009876', '天弘甄选食品饮料股票C', '1.89', '84.12', '4.42', '0.0835', 8),
    @('009875', '天弘甄选食品饮料股票A', '1.59', '84.12', '4.42', '0.0703', 8),
    @('001030', '天弘云端生活优选灵活配置混合A', '1.40', '79.61', '4.26', '0.0596', 6),
    @('009954', '北信瑞丰优选成长股票', '0.55', '93.35', '4.35', '0.0239', 10),
    @('200016', '长城稳健成长灵活配置混合', '0.73', '79.86', '3.01', '0.0220', 5),
    @('010331', '天弘消费股票A', '0.28', '85.31', '3.91', '0.0109', 8),
    @('010332', '天弘消费股票C', '0.21', '85.31', '3.91', '0.0082', 8),
    @('005015', '泰康景泰回报混合C', '0.39', '34.29', '1.50', '0.0058', 5),
    @('015462', '天弘云端生活优选灵活配置混合C', '0.00', '79.61', '4.26', 'NUM:0', 6)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    # A: 0-based row index, same bordered/bold style used throughout
    $aCell = $q3.Cells.Item($r, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    # B: fund code -- force text so leading zeros survive
    $bCell = $q3.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[0]

    # C: fund name (plain text)
    $q3.Cells.Item($r, 3).Value = $row[1]

    # D, E, F: size / position / weight -- stored as text in the source data
    $dCell = $q3.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[2]

    $eCell = $q3.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[3]

    $fCell = $q3.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[4]

    # G: market value (亿元) -- text, except the literal-zero special case
    $gRaw = $row[5]
    $gCell = $q3.Cells.Item($r, 7)
    if ($gRaw.StartsWith("NUM:")) {
        $gCell.Value = [double]($gRaw.Substring(4))
    } else {
        $gCell.NumberFormat = "@"
        $gCell.Value = $gRaw
    }

    # H: position rank (number)
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Prepend the 2022-Q3 summary row to "总计"
# ---------------------------------------------------------------------------

$summaryDates = @("2022-Q3", "2022-Q2", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$summaryCounts = @(15, 29, 37, 44, 9, 32, 39, 30)
$summaryValues = @(15.36, 20.56, 22.04, 27.95, 1.54, 28.8, 30.32, 22.27)

# Make sure the newly-appended row 9 carries the same styling as the other
# data rows (copy column A's formatting down from row 2).
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $summaryDates.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $summaryDates[$i]
    $totalSheet.Cells.Item($r, 3).Value = $summaryCounts[$i]
    $totalSheet.Cells.Item($r, 4).Value = $summaryValues[$i]
}
